# Update xGA (away) values for Lens (row 6) for Giornata 11-18,
# reflecting the recalculated Ligue1 output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K6").Value = 4.492175
$ws.Range("L6").Value = 4.492175
$ws.Range("M6").Value = 5.164258
$ws.Range("N6").Value = 5.164258
$ws.Range("O6").Value = 5.66808
$ws.Range("P6").Value = 5.66808
$ws.Range("Q6").Value = 6.596571
$ws.Range("R6").Value = 7.756951
$ws.Range("S6").Value = 7.756951
